# Error Calculations and Plots
# Re-derive the "missing data" mask: drop two rows entirely (RM 232, SC 92)
# and re-randomize which cells in columns C:F ("B","C","D","F") are blanked
# out vs populated for several remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove rows for "RM 232" (row 26) and "SC 92" (row 28) --------------
# Delete the lower row first so the earlier row index is unaffected.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- After the deletions the remaining rows shifted up; re-apply the new
#     missing/non-missing pattern for columns C:F on the affected rows ----
$ws.Range("F5").Value = $null
$ws.Range("F11").Value = 17.65
$ws.Range("E19").Value = -6.5
$ws.Range("F19").Value = $null
$ws.Range("E21").Value = $null
$ws.Range("E23").Value = -7
$ws.Range("F23").Value = 16.48
$ws.Range("F25").Value = 16.6
$ws.Range("C26").Value = $null
$ws.Range("C27").Value = 10
$ws.Range("E27").Value = $null
$ws.Range("F27").Value = $null
$ws.Range("C29").Value = $null
$ws.Range("F29").Value = $null
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
